# sprint3 final DOCs update
# Burn-down tracker: the team re-estimated the total backlog for Sprint 3
# from 37 points up to 45 points. Everything else on Sheet1 (the daily
# "BurnDn" / "Ideal" columns, the backlog total in C26, and the chart
# series that plot them) is driven off this single input cell by
# formulas already in the sheet, so updating C2 is the whole edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated total backlog for the sprint.
$ws.Range("C2").Value = 45

# Leave the cursor where the author left it when they saved.
$ws.Range("D4").Select()
